$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 27.2
$ws.Range("B3").Value = 34.7
$ws.Range("C3").Value = 30.5
$ws.Range("C4").Value = 32.8
$ws.Range("C9").Value = 34.5
$ws.Range("C10").Value = 33.8
$ws.Range("C13").Value = 28.5
$ws.Range("C21").Value = 27.1
$ws.Range("C22").Value = 31
